$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.054352781976451929
$ws.Range("B1").Value = 0.054352781942259079

$ws.Range("A2").Value = 0.042095558335632527
$ws.Range("B2").Value = -0.042095558407180184

$ws.Range("A3").Value = -0.050969206730113945
$ws.Range("B3").Value = 0.050969206695782705
